$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 3 (the "abc@gmail.com" / "Pass1234" row), which shifts
# the row below it (solih48913@gamebcs.com) up to become the new row 3.
$ws.Rows.Item(3).Delete()

# Excel leaves the selection on the row that now occupies the deleted row's
# position, with the whole row selected.
[void]$ws.Rows.Item(3).Select()
